# Update the "Förändrad" (changed) date column (C) for rows 2-6
# from serial date 45208 (2023-10-09) to 45212 (2023-10-13),
# keeping the existing date formatting/style on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45208) {
        $cell.Value = 45212
    }
}
